$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column H header "Save", matching the formatting (bold/border) used by G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"
$excel.CutCopyMode = $false

# Values for H2:H20 ("Save" flag) - default 0, with 1 for specific rows
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
